$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -0.1905
$ws.Range("C6").Value = -0.2534
$ws.Range("C14").Value = -0.0962
$ws.Range("C15").Value = -0.1881
$ws.Range("C16").Value = 0.1578
$ws.Range("C17").Value = 0.0071
$ws.Range("C19").Value = 0.2071
$ws.Range("C20").Value = -0.3541
